# Fix the "property_category" column values on the 建物 (Building) and
# 汽車 (Car) sheets. Both sheets had every row incorrectly tagged with the
# "land" category (copy/paste leftover); this corrects them to "building"
# and "car" respectively. See commit "#5: property aircraft done".

$wb = $excel.ActiveWorkbook

# 建物 (Building) sheet: column I is "property_category", data rows 2-7
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2:I7").Value = "building"

# 汽車 (Car) sheet: column H is "property_category", data row 2
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2:H2").Value = "car"
